$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-20 down to 12-21
$ws.Rows("11:11").Insert()

# Populate the new row 11 with data
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44554
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112043
$ws.Cells.Item(11, 7).Value = "Pepino dulce"
$ws.Cells.Item(11, 8).Value = "Cultivar XV región"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(11, 11).Value = 5000
$ws.Cells.Item(11, 12).Value = 6000
$ws.Cells.Item(11, 13).Value = 5500
$ws.Cells.Item(11, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 550
$ws.Cells.Item(11, 17).Value = 10
$ws.Cells.Item(11, 18).Value = "Hortaliza"
